$d = $word.ActiveDocument

$replacements = @(
    @{old = "831÷9="; new = "374÷8="},
    @{old = "869÷9="; new = "251÷8="},
    @{old = "775÷9="; new = "643÷6="},
    @{old = "824÷3="; new = "844÷8="},
    @{old = "494÷6="; new = "302÷9="},
    @{old = "532÷7="; new = "812÷7="},
    @{old = "627÷7="; new = "232÷6="},
    @{old = "501÷5="; new = "852÷6="},
    @{old = "942÷8="; new = "763÷5="},
    @{old = "969÷3="; new = "871÷4="},
    @{old = "901÷7="; new = "135÷6="},
    @{old = "576÷4="; new = "252÷2="},
    @{old = "151÷7="; new = "697÷5="},
    @{old = "674÷4="; new = "660÷6="},
    @{old = "985÷3="; new = "245÷2="},
    @{old = "805÷8="; new = "834÷7="},
    @{old = "487÷8="; new = "905÷6="},
    @{old = "651÷8="; new = "373÷8="},
    @{old = "557÷8="; new = "392÷7="},
    @{old = "990÷6="; new = "307÷5="},
    @{old = "263÷8="; new = "625÷7="},
    @{old = "691÷8="; new = "153÷2="},
    @{old = "950÷9="; new = "167÷8="},
    @{old = "378÷4="; new = "131÷5="},
    @{old = "634÷2="; new = "293÷4="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
